# Apply the "representative" (next-of-kin contact) block to the
# UserCreation sheet, add the matching hyperlink for the representative's
# email, widen column A, and move the sheet view/selection down to the
# newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserCreation")

# --- New key/value rows 43-54 -------------------------------------------
$ws.Range("A43").Value = "representativeTitle"
$ws.Range("B43").Value = "Mr"

$ws.Range("A44").Value = "representativeFirstNam"
$ws.Range("B44").Value = "repAutomationFirstName"

$ws.Range("A45").Value = "representativeMiddleName"
$ws.Range("B45").Value = "repAutomationMiddleName"

$ws.Range("A46").Value = "representativeLastName"
$ws.Range("B46").Value = "repAutomationLastName"

$ws.Range("A47").Value = "relationship"
$ws.Range("B47").Value = "Cousin"

$ws.Range("A48").Value = "representativeAddLine1"
$ws.Range("B48").Value = "Mrs Smith 98 Shirley Street"

$ws.Range("A49").Value = "representativeAddLine2"
$ws.Range("B49").Value = "PIMPAMA QLD 4209"

$ws.Range("A50").Value = "representativeSetSuburb"
$ws.Range("B50").Value = "Aarons Pass"
$ws.Range("C50").Value = "Arrawarra Headland, New South Wales, 2456"

$ws.Range("A51").Value = "representativePhone"
$ws.Range("B51").Value = "(08) 4356-7689"

$ws.Range("A52").Value = "representativeMobile"
$ws.Range("B52").Value = "(08) 4356-7689"

$ws.Range("A53").Value = "representativeEmail"
$ws.Range("B53").Value = "jtariq@ucm.com.au"

$ws.Range("A54").Value = "representativeDrivingLicenseNo"
$ws.Range("B54").Value = "PIMPAMA QLD 4209"

# --- Hyperlink on the new representative email cell ----------------------
$ws.Hyperlinks.Add($ws.Range("B53"), "mailto:jtariq@ucm.com.au")

# --- Column A got a bit wider to fit the new "representative..." keys ---
$ws.Columns("A").ColumnWidth = 33 + 1/6

# --- Scroll the view down to the newly added rows and reselect ----------
$ws.Activate()
$ws.Range("B50").Select()

Write-Host "UserCreation sheet updated with representative fields"
